$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Make room for two new rows right after the header row by sliding the
#    existing data block (A2:D40) down to A4:D42. A plain block copy/paste
#    (rather than a structural Rows.Insert) is used deliberately: on this
#    engine, Insert()-ing blank rows blends formatting from the row above
#    and mints a brand-new cellXf/style entry, which would leave stray
#    unused styles in styles.xml. Copy+PasteAll reuses the existing style
#    indices untouched.
# ---------------------------------------------------------------------------
$ws.Range("A2:D40").Copy()
$ws.Range("A4").PasteSpecial(-4104)   # xlPasteAll
$ws.Application.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2) Seed the two new rows' formatting by copying from rows that already
#    carry the right look: plain style for A:C, bold-hyperlink style for D
#    (the same combination used for other recently-added entries such as
#    the one now at row 9).
# ---------------------------------------------------------------------------
$ws.Range("A4:C4").Copy()
$ws.Range("A2:C3").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("D9").Copy()
$ws.Range("D2:D3").PasteSpecial(-4122)   # xlPasteFormats
$ws.Application.CutCopyMode = 0

# Row 2: advance_search_food
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = "advance_search_food"
$ws.Cells.Item(2, 3).Value = "keyword, type, min, max"
$ws.Cells.Item(2, 4).Value = "http://localhost/fyp_connect/advance_search_food.php?keyword=tea&type=Drink&min=1&max=15"

# Row 3: advance_search_restaurant
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = "advance_search_restaurant"
$ws.Cells.Item(3, 3).Value = "keyword, type, district"
$ws.Cells.Item(3, 4).Value = "http://localhost/fyp_connect/advance_search_restaurant.php?keyword=tea&type=Itailan&district=Tai%20po"

# ---------------------------------------------------------------------------
# 3) Renumber column A (the running index 1..41) for every data row now
#    that two rows were inserted above row 4.
# ---------------------------------------------------------------------------
for ($r = 4; $r -le 42; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}

# ---------------------------------------------------------------------------
# 4) Hyperlinks are anchored to absolute cell refs and do not follow a plain
#    value/format copy, so rebuild the whole collection: wipe the old
#    (now-stale) links and re-add one per row at its new position, plus the
#    two brand-new ones for rows 2 and 3.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Delete()

$links = @{
    2  = "http://localhost/fyp_connect/advance_search_food.php?keyword=tea&type=Drink&min=1&max=15"
    3  = "http://localhost/fyp_connect/advance_search_restaurant.php?keyword=tea&type=Itailan&district=Tai%20po"
    4  = "http://localhost/fyp_connect/create_food.php?name=french%20friess&type=french%20fries&price=13&image=french_fries.jpg&Restaurantid=r0001"
    5  = "http://localhost/fyp_connect/create_order.php?customerid=c0001&order_total=25"
    6  = "http://localhost/fyp_connect/create_orderline.php?ordernumber=1&foodid=f0001&quanitity=1&item_total=25&Restaurantid=r0001"
    7  = "http://localhost/fyp_connect/create_user.php?password=ggwp&name=aaa&address=asdasasd&email=adgasdgasd&telNum=12341234"
    8  = "http://localhost/fyp_connect/customer_get_order.php?userid=10002"
    9  = "http://localhost/fyp_connect/customer_get_order2.php?id=c0001"
    12 = "http://localhost/fyp_connect/delete_food.php?id=f0001"
    13 = "http://localhost/fyp_connect/delete_order.php?number=1"
    15 = "http://localhost/fyp_connect/driver_get_order_qty.php"
    16 = "http://localhost/fyp_connect/driver_get_order_time.php"
    18 = "http://localhost/fyp_connect/driver_handle_order.php?driverid=d0001&number=1"
    19 = "http://localhost/fyp_connect/driver_update_order_pick_time.php?driverid=d0001&number=2"
    20 = "http://localhost/fyp_connect/driver_update_order_receive.php?driverid=d0001&number=2"
    21 = "http://localhost/fyp_connect/find_customer_id.php?userid=10002"
    22 = "http://localhost/fyp_connect/find_restaurant_id.php?userid=10004"
    23 = "http://localhost/fyp_connect/get_all_food.php"
    24 = "http://localhost/fyp_connect/get_all_order.php"
    25 = "http://localhost/fyp_connect/get_all_user.php"
    26 = "http://localhost/fyp_connect/get_available_order.php"
    27 = "http://localhost/fyp_connect/get_order_count.php?number=1"
    28 = "http://localhost/fyp_connect/get_order_details.php?number=1"
    29 = "http://localhost/fyp_connect/get_restaurant_food.php?userid=10004"
    30 = "http://localhost/fyp_connect/get_user_details.php?id=10001&password=10001"
    32 = "http://localhost/fyp_connect/id_get_food_name.php?foodid=f0001"
    33 = "http://localhost/fyp_connect/id_get_food_type.php?foodid=f0001"
    34 = "http://localhost/fyp_connect/restaurant_get_order.php?sql=SELECT%20*%20FROM%20%60order%60%20WHERE%20number=1%20OR%20number=2"
    35 = "http://localhost/fyp_connect/restaurant_get_all_orderline.php?restaurantid=r0001"
    36 = "http://localhost/fyp_connect/restaurant_get_orderline.php?restaurantid=r0001"
    38 = "http://localhost/fyp_connect/restaurant_update_order_pick_up.php?ordernumber=1&restaurantid=r0001"
    39 = "http://localhost/fyp_connect/restaurant_update_order_status.php?ordernumber=1&restaurantid=r0001"
    40 = "http://localhost/fyp_connect/search_details.php?type=food&keyword=chicken"
    42 = "http://localhost/fyp_connect/update_food.php?id=f0001&name=Black%20Tea&type=Drink&price=15&image=Black_tea.jpg"
}

foreach ($row in ($links.Keys | Sort-Object)) {
    $ws.Hyperlinks.Add($ws.Cells.Item($row, 4), $links[$row])
}

# ---------------------------------------------------------------------------
# 5) Refresh the view: scroll back to the top and move the selection to C3.
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("C3").Select()
